# Update "Recommandations" sheet (sheet1): rows 2-43 get new values
# (re-ranked table with reordered/renamed rows and refreshed numbers),
# and the 4 trailing rows (44-47) are removed entirely.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

$sheet1Data = @(
    @('BRVM - SERVICES PUBLICS', 0, 6, 2414.05, 100.18, '🟡 Observer', '➖ Neutre'),
    @('CFAO MOTORS CI', 0, 3, 1945, 650, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AUTRES SECTEURS', 0, 3, 1892.68, 640.71, '🟡 Observer', '➖ Neutre'),
    @('NEI-CEDA CI', 0, 3, 1785, 595, '🟡 Observer', '➖ Neutre'),
    @('UNIWAX CI', 0, 3, 1740, 585, '🟡 Observer', '➖ Neutre'),
    @('SETAO CI', 0, 3, 1720, 560, '🟡 Observer', '➖ Neutre'),
    @('AIR LIQUIDE CI', 0, 3, 1610, 540, '🟡 Observer', '➖ Neutre'),
    @('BRVM - DISTRIBUTION', 0, 3, 1095.2, 367.96, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TRANSPORT', 0, 3, 1057.41, 352.47, '🟡 Observer', '➖ Neutre'),
    @('SAFCA CI', 0, 1, 985, 985, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AGRICULTURE', 0, 3, 964.38, 323.48, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE', 0, 3, 787.58, 262.99, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE', 0, 3, 650.7, 217.65, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL', 0, 3, 565.31, 189.11, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIELS', 0, 3, 413.84, 137.32, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRESTIGE', 0, 3, 385.85, 127.91, '🟡 Observer', '➖ Neutre'),
    @('BRVM - FINANCES', 0, 3, 363.34, 120.64, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES FINANCIERS', 0, 3, 357.08, 118.56, '🟡 Observer', '➖ Neutre'),
    @('BRVM - ENERGIE', 0, 3, 325.49, 109.47, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 3, 312.28, 104.87, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TELECOMMUNICATIONS', 0, 3, 275.49, 92.11, '🟡 Observer', '➖ Neutre'),
    @('SAFCA CI (SAFC)', 3, 0, 21.29, 7.11, '🟢 Achat', '✅ Renforcer'),
    @('BERNABE CI (BNBC)', 2, 0, 9.39, 2.19, '🟡 Observer', '➖ Neutre'),
    @('VIVO ENERGY CI (SHEC)', 2, 0, 8.44, 5.09, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA ML (BOAM)', 1, 0, 7.37, 7.37, '🟡 Observer', '➖ Neutre'),
    @('SAPH CI (SPHC)', 1, 0, 4.69, 4.69, '🟡 Observer', '➖ Neutre'),
    @('SERVAIR ABIDJAN CI (ABJC)', 1, 0, 4.07, 4.07, '🟡 Observer', '➖ Neutre'),
    @('SICABLE CI (CABC)', 1, 0, 4, 4, '🟡 Observer', '➖ Neutre'),
    @('AIR LIQUIDE CI (SIVC)', 1, 1, 2.94, -1.82, '🟡 Observer', '👀 À surveiller'),
    @('ORAGROUP TOGO (ORGT)', 1, 1, 0.96, 2.53, '🟡 Observer', '👀 À surveiller'),
    @('ECOBANK TRANS. INCORP. TG (ETIT)', 1, 1, 0.42, 6.67, '🟡 Observer', '👀 À surveiller'),
    @('TOTAL', 0, 3, 0, 0, '🟡 Observer', '➖ Neutre'),
    @('ORANGE COTE D''IVOIRE (ORAC)', 0, 1, -1.06, -1.06, '🟡 Observer', '➖ Neutre'),
    @('ECOBANK COTE D''''IVOIRE (ECOC)', 0, 1, -1.32, -1.32, '🟡 Observer', '➖ Neutre'),
    @('SONATEL SN (SNTS)', 0, 1, -1.4, -1.4, '🟡 Observer', '➖ Neutre'),
    @('LOTERIE NATIONALE DU BENIN (LNBB)', 1, 1, -1.58, 5.88, '🟡 Observer', '👀 À surveiller'),
    @('FILTISAC CI (FTSC)', 0, 1, -1.84, -1.84, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA NG (BOAN)', 0, 1, -2.71, -2.71, '🟡 Observer', '➖ Neutre'),
    @('SETAO CI (STAC)', 0, 1, -3.45, -3.45, '🟡 Observer', '➖ Neutre'),
    @('SOCIETE GENERALE COTE D''IVOIRE (SGBC)', 0, 1, -3.71, -3.71, '🟡 Observer', '➖ Neutre'),
    @('CIE CI (CIEC)', 0, 2, -4.98, -2.27, '🟡 Observer', '➖ Neutre'),
    @('UNIWAX CI (UNXC)', 0, 1, -5.17, -5.17, '🟡 Observer', '➖ Neutre')
)

$rows1 = $sheet1Data.Count
$cols1 = 7
$arr1 = New-Object 'object[,]' $rows1,$cols1
for ($i = 0; $i -lt $rows1; $i++) {
    for ($j = 0; $j -lt $cols1; $j++) {
        $arr1[$i,$j] = $sheet1Data[$i][$j]
    }
}
$ws1.Range("A2:G43").Value = $arr1

# The source table shrank from 47 to 43 data+header rows; remove the now
# obsolete trailing rows so the sheet's dimension becomes A1:G43.
$ws1.Rows("44:47").Delete()

# Update "Top_YTD" sheet (sheet2): only column B (Progression YTD %) changes.
$ws2 = $wb.Worksheets.Item("Top_YTD")

$sheet2Data = @(
    @('BRVM - SERVICES PUBLICS', 420906.54),
    @('CFAO MOTORS CI', 41806.25),
    @('BRVM - AUTRES SECTEURS', 38921.44),
    @('NEI-CEDA CI', 33470.24),
    @('UNIWAX CI', 31341.5),
    @('SETAO CI', 30418.4),
    @('AIR LIQUIDE CI', 25700),
    @('BRVM - DISTRIBUTION', 9958.45),
    @('BRVM - TRANSPORT', 9163.38),
    @('BRVM - AGRICULTURE', 7386.16)
)

$rows2 = $sheet2Data.Count
$arr2 = New-Object 'object[,]' $rows2,1
for ($i = 0; $i -lt $rows2; $i++) {
    $arr2[$i,0] = $sheet2Data[$i][1]
}
$ws2.Range("B2:B11").Value = $arr2
